# Repull data, push all data, mean calculation
# Update column F (dSF) values for several rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = 2
    13 = 2
    15 = 0
    21 = 2
    25 = 2
    26 = 1
    28 = -1
    33 = 1
    34 = -1
    35 = 0
    39 = 0
    42 = 1
    45 = 2
    50 = -1
    54 = -1
    60 = 2
    61 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
